$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("이준기")

# Fill in row 22 with the new time-log entry
$ws.Range("A22").Value = "10월 29일"
$ws.Range("B22").Value = 0.5
$ws.Range("C22").Value = 0.55208333333333337
$ws.Range("D22").Value = 0
$ws.Range("E22").Value = 75
$ws.Range("F22").Value = "Activity Diagram 전처리 파트 수정"

# Update the active selection to reflect the latest edit position
$ws.Range("G25").Select()
